# Edit workbook per requested change:
#  - within every year block of 4 quarter rows (A,B,C,D), swap the
#    contents of the "B" quarter row and the "C" quarter row
#    (label in col A together with the B:E data), leaving the "A" and
#    "D" quarter rows untouched
#  - remove columns F (家用洗衣机产销率) and G (家用洗衣机销售量)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers of the "B" / "C" quarter rows that need to be swapped,
# given as consecutive pairs (Brow, Crow).
$pairs = @(3,4, 7,8, 11,12, 15,16, 19,20, 23,24, 27,28, 31,32, 35,36, 39,40, 43,44, 47,48, 51,52, 55,56, 59,60, 63,64, 67,68, 71,72, 75,76, 79,80)

# scratch row used as a temporary holding area while swapping two rows;
# it is well outside the populated data range (A1:G81)
$tmpRow = 200

for ($i = 0; $i -lt $pairs.Length; $i += 2) {
    $r1 = $pairs[$i]
    $r2 = $pairs[$i + 1]

    $rng1 = $ws.Range("A" + $r1 + ":E" + $r1)
    $rng2 = $ws.Range("A" + $r2 + ":E" + $r2)
    $tmp  = $ws.Range("A" + $tmpRow + ":E" + $tmpRow)

    # Copy() onto a non-blank destination does not clear cells that the
    # source has blank, so every destination must be cleared immediately
    # before it receives its new content.
    $tmp.ClearContents()
    $rng1.Copy($tmp)

    $rng1.ClearContents()
    $rng2.Copy($rng1)

    $rng2.ClearContents()
    $tmp.Copy($rng2)

    $ws.Rows($tmpRow).Delete()
}

# Drop the now-unneeded F (产销率) and G (销售量) columns entirely.
$ws.Columns("F:G").Delete()
